$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.206311106681824
$ws.Range("B1").Value = 1.367578744888306
$ws.Range("C1").Value = 6.955246925354004
$ws.Range("D1").Value = 2.187525272369385
$ws.Range("E1").Value = 1.168618202209473
